$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.655.95"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.426.63"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "552.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.71"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +1.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.159"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.78%  "
$ws.Range("E10").Value = "  -0.66%  "
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").Value = "67.598.28"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.35"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "330.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("E18").Value = "  -1.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -0.98%  "
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "416.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.70%  "
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.27%  "
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("E33").Value = "  +0.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.105"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.58%  "
$ws.Range("E36").Value = "  -2.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.25"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.34%  "
$ws.Range("E39").Value = "  -2.37%  "
$ws.Range("E40").Value = "  -2.45%  "
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "129.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0706"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.82%  "
$ws.Range("E44").Value = "  -0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.553"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0914"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("E48").Value = "  -6.89%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.00%  "
$ws.Range("E50").Value = "  +2.37%  "
$ws.Range("E51").Value = "  +0.11%  "
